$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2404")
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2410")
}

# Convert the data range into an Excel Table (ListObject)
$range = $ws.Range("A1:U55")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# Freeze the top row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
